$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build formatting (bold font, thin box border, center/top alignment) on B1
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160
$b1.Borders.LineStyle = 1
$b1.Borders.Weight = 2

# Copy the resulting format onto A2 (keeps style table minimal / avoids
# generating extra intermediate cellXfs entries)
$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
